$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The parser table used to have sample rows 4-8 populated; clear their
# contents (keep formatting) now that the real parser populates them.
$ws.Range("A4:G8").ClearContents()

# Leave the selection on the first cleared block, matching the state the
# workbook was left in after the cleanup.
$ws.Range("A4:G5").Select()
